$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7540575265884399
$ws.Range("B1").Value = 0.7489599585533142
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.515378832817078
$ws.Range("E1").Value = 0.9225829243659973
